# Add a new "Greece" sheet, modeled on the existing "Croatia" sheet,
# with market/part values updated for Greece, and make it the active tab.

$wb = $excel.ActiveWorkbook

$croatia = $wb.Worksheets.Item("Croatia")

# Select the whole sheet on Croatia before copying - this mirrors what
# happened to the source sheet's view state when the new sheet was created.
$croatia.Activate() | Out-Null
$croatia.Cells.Select() | Out-Null

# Duplicate Croatia into a new sheet placed right after it.
$croatia.Copy($null, $croatia)

# The copy becomes the last sheet in the workbook.
$greece = $wb.Worksheets.Item($wb.Worksheets.Count)
$greece.Name = "Greece"

# Update the market name and part number for the Greece sheet.
$greece.Range("B2").Value = "Greece Market"
$greece.Range("B4").Value = "NGC-4119/T3190"

# Make the new sheet the active one, with B4 selected.
$greece.Activate() | Out-Null
$greece.Range("B4").Select() | Out-Null
